$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear "DATA CADASTRO" (column O) for all data rows (2 through 255)
$ws.Range("O2:O255").ClearContents()

# Clear "DATA FECHAMENTO" (column P) for row 179 only
$ws.Range("P179").ClearContents()
